$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Entrenadores")
$ws.Select()
